$d = $word.ActiveDocument

function Replace-ParagraphXml {
    param($FindText, $InnerXml)

    $rng = $d.Content
    $found = $rng.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw ("Replace-ParagraphXml: text not found: " + $FindText)
    }
    $rng.Delete()

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' + $InnerXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData>' +
        '</pkg:part>' +
        '</pkg:package>'

    $rng.InsertXML($xml)
}

# 1. "recent history related to the idea or product development. Keep this to around a paragraph."
Replace-ParagraphXml "recent history related to the idea or product development. Keep this to around a paragraph." (
    '<w:p w14:paraId="5AF8EA69" w14:textId="77777777" w:rsidR="00976A67" w:rsidRPr="00976A67" w:rsidRDefault="00976A67" w:rsidP="00976A67">' +
      '<w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' +
      '<w:proofErr w:type="gramStart"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>recent</w:t></w:r>' +
      '<w:proofErr w:type="gramEnd"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> history related to the idea or product development. Keep this to around a paragraph.</w:t></w:r>' +
    '</w:p>'
)

# 2. "one sentence. If you can't keep it to one sentence, try to create a short list of goals (no more"
Replace-ParagraphXml "one sentence. If you can’t keep it to one sentence, try to create a short list of goals (no more" (
    '<w:p w14:paraId="7C825D98" w14:textId="77777777" w:rsidR="00976A67" w:rsidRPr="00976A67" w:rsidRDefault="00976A67" w:rsidP="00976A67">' +
      '<w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' +
      '<w:proofErr w:type="gramStart"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>one</w:t></w:r>' +
      '<w:proofErr w:type="gramEnd"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> sentence. If you can’t keep it to one sentence, try to create a short list of goals (no more</w:t></w:r>' +
    '</w:p>'
)

# 3. "than 4!)"
Replace-ParagraphXml "than 4!)" (
    '<w:p w14:paraId="304D7DFD" w14:textId="77777777" w:rsidR="00976A67" w:rsidRPr="00976A67" w:rsidRDefault="00976A67" w:rsidP="00976A67">' +
      '<w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' +
      '<w:proofErr w:type="gramStart"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>than</w:t></w:r>' +
      '<w:proofErr w:type="gramEnd"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> 4!)</w:t></w:r>' +
    '</w:p>'
)

# 4. "long and where."
Replace-ParagraphXml "long and where." (
    '<w:p w14:paraId="4EEEA52A" w14:textId="77777777" w:rsidR="00976A67" w:rsidRPr="00976A67" w:rsidRDefault="00976A67" w:rsidP="00976A67">' +
      '<w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' +
      '<w:proofErr w:type="gramStart"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>long</w:t></w:r>' +
      '<w:proofErr w:type="gramEnd"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> and where.</w:t></w:r>' +
    '</w:p>'
)

# 5. Schedule paragraph: split off "Live" with proofErr around it.
Replace-ParagraphXml "Food Provider and Food Consumer onboarding – 1/14 – 3/1, Live Deployment – 3/1" (
    '<w:p w14:paraId="0BD91896" w14:textId="77777777" w:rsidR="00976A67" w:rsidRPr="00976A67" w:rsidRDefault="00976A67" w:rsidP="00976A67">' +
      '<w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' +
      '<w:r w:rsidRPr="00976A67"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">9. </w:t></w:r>' +
      '<w:r w:rsidRPr="00976A67"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:b/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Schedule</w:t></w:r>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">: Build Schedule – 1/30 – 2/13, Food Provider and Food Consumer onboarding – 1/14 – 3/1, </w:t></w:r>' +
      '<w:proofErr w:type="gramStart"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Live</w:t></w:r>' +
      '<w:proofErr w:type="gramEnd"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> Deployment – 3/1</w:t></w:r>' +
    '</w:p>'
)

# 6. Doug runs paragraph: split off "it to." with proofErr around it.
Replace-ParagraphXml "Doug would like to donate this food but doesn’t know where to send it to. He also doesn’t want to spend extra time or effort delivering it somewhere." (
    '<w:p w14:paraId="2C9455A1" w14:textId="77777777" w:rsidR="007C7AB8" w:rsidRDefault="007C7AB8" w:rsidP="00976A67">' +
      '<w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t xml:space="preserve">Doug runs a catering business and caters to large events on a weekly basis. Doug often finds himself with plenty of leftover food after an event that usually gets thrown out. Doug would like to donate this food but doesn’t know where to send </w:t></w:r>' +
      '<w:proofErr w:type="gramStart"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t>it to.</w:t></w:r>' +
      '<w:proofErr w:type="gramEnd"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t xml:space="preserve"> He also doesn’t want to spend extra time or effort delivering it somewhere.</w:t></w:r>' +
    '</w:p>'
)

# 7. "For Food Banks within a 10 mile radius of you -what donation method would you be willing to do?"
Replace-ParagraphXml "For Food Banks within a 10 mile radius of you -what donation method would you be willing to do?" (
    '<w:p w14:paraId="3E529291" w14:textId="61F2198C" w:rsidR="004F2AE4" w:rsidRDefault="004F2AE4" w:rsidP="004F2AE4">' +
      '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t xml:space="preserve">For Food Banks within a </w:t></w:r>' +
      '<w:proofErr w:type="gramStart"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t>10 mile</w:t></w:r>' +
      '<w:proofErr w:type="gramEnd"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t xml:space="preserve"> radius of you -what donation method would you be willing to do?</w:t></w:r>' +
    '</w:p>'
)

# 8. "How do you typically measure leftover food? (ie: servings, weight etc.)"
Replace-ParagraphXml "How do you typically measure leftover food? (ie: servings, weight etc.)" (
    '<w:p w14:paraId="45E0A874" w14:textId="545EEAB3" w:rsidR="004F2AE4" w:rsidRDefault="004F2AE4" w:rsidP="004F2AE4">' +
      '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t>How do you typically measure leftover food? (</w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:proofErr w:type="gramStart"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t>ie</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:proofErr w:type="gramEnd"/>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t>: servings, weight etc.)</w:t></w:r>' +
    '</w:p>'
)

# 9. Merge the 6 "For Vendors within that radius..." runs into a single run.
Replace-ParagraphXml "For Vendors within that radius - what donation method would you be willing to do?" (
    '<w:p w14:paraId="021B331B" w14:textId="77777777" w:rsidR="004F2AE4" w:rsidRDefault="004F2AE4" w:rsidP="004F2AE4">' +
      '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t>For Vendors within that radius - what donation method would you be willing to do?</w:t></w:r>' +
    '</w:p>'
)

# 10 & 11. Add "?" run to "What information..." paragraph, then insert a new paragraph after it
#          (moving the _GoBack bookmark to the end of the new paragraph).
Replace-ParagraphXml "What information about the food would you like to see" (
    '<w:p w14:paraId="493AF492" w14:textId="54E830AC" w:rsidR="004F2AE4" w:rsidRPr="004F2AE4" w:rsidRDefault="004F2AE4" w:rsidP="004F2AE4">' +
      '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t>What information about the food would you like to see</w:t></w:r>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t>?</w:t></w:r>' +
    '</w:p>' +
    '<w:p>' +
      '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/></w:rPr><w:t>What credentials do you have that show you are a 501c/non-profit?</w:t></w:r>' +
      '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
      '<w:bookmarkEnd w:id="0"/>' +
    '</w:p>'
)
